$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# This CI run re-generated the handback report: the file that was already
# tracked (old GUID) got a fresh pass (new timestamps / new handed-back
# xliff hashes) and a second file newly showed up as row 3 on every sheet.
# ---------------------------------------------------------------------------

$oldGuid = "7925027b-ecc4-4ad7-80b0-2392b7823ebf"
$oldMd   = "$oldGuid.md"

$guid1   = "0a56a075-b996-4dae-82d4-58a1bbedd02c"
$md1     = "$guid1.md"
$md1Disp = "e2e\$md1"
$hash1   = "58f9cdb8358c1945ca430bc0ec0939f10a56694e"
$zhXlf1  = "$guid1.$hash1.zh-cn.xlf"
$deXlf1  = "$guid1.$hash1.de-de.xlf"

$guid2   = "0f6864f6-bed6-4170-87d7-20c0f4e5f3d6"
$md2     = "$guid2.md"
$md2Disp = "e2e\$md2"
$hash2   = "f40a51c048fab9b0d50fee9ac333d5286fe99af4"
$zhXlf2  = "$guid2.$hash2.zh-cn.xlf"
$deXlf2  = "$guid2.$hash2.de-de.xlf"

$status    = "Handed back: in sync with en-US"
$extDot    = ".md"
$extNoDot  = "e2e"
$htStr     = "ht"
$falseStr  = "False"
$trueStr   = "True"

$genDate1        = "2016-08-12 03:22:30"
$handoffDate1    = "2016-08-12 03:22:25"
$handbackZhDate1 = "2016-08-12 03:22:42"
$handbackDeDate1 = "2016-08-12 03:22:49"

$genDate2        = "2016-08-12 03:22:30"
$handoffDate2    = "2016-08-12 03:22:25"
$handbackZhDate2 = "2016-08-12 03:22:42"
$handbackDeDate2 = "2016-08-12 03:22:49"

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# Colour used by the workbook's custom "HyperLink" cell style (RGB 6495ED,
# stored as BGR for the COM Color property).
$linkColor = 15570276

function Style-AsHyperlink($rng) {
    $rng.Font.Underline = $true
    $rng.Font.Color = $linkColor
}

function Set-Hyperlink($ws, $addr, $url, $display) {
    $rng = $ws.Range($addr)
    if ($rng.Hyperlinks.Count -gt 0) {
        $rng.Hyperlinks.Delete()
    }
    $ws.Hyperlinks.Add($rng, $url, "", "", $display) | Out-Null
    Style-AsHyperlink $rng
}

# ---------------------------------------------------------------------------
# Sheet "Overview" - one summary row per handled file.
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Row 2: the already-known file, re-processed with new timestamp.
$ov.Range("A2").Value = $md1
$ov.Range("B2").Value = $md1Disp
$ov.Range("G2").Value = $genDate1
$ov.Range("G2").NumberFormat = $dateFmt
Set-Hyperlink $ov "B2" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md1" $md1Disp

# Row 3: newly appeared file.
$ov.Range("A3").Value = $md2
$ov.Range("B3").Value = $md2Disp
$ov.Range("C3").Value = $extDot
$ov.Range("E3").Value = $status
$ov.Range("F3").Value = $status
$ov.Range("G3").Value = $genDate2
$ov.Range("G3").NumberFormat = $dateFmt
Set-Hyperlink $ov "B3" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md2" $md2Disp

$ovTable = $ov.ListObjects.Item(1)
$ovTable.Resize($ov.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

# Row 2: re-processed file - new xliff hash + new handoff/handback timestamps.
$zh.Range("A2").Value = $md1
$zh.Range("G2").Value = $zhXlf1
$zh.Range("H2").Value = $handoffDate1
$zh.Range("H2").NumberFormat = $dateFmt
$zh.Range("I2").Value = $md1
$zh.Range("J2").Value = $zhXlf1
$zh.Range("K2").Value = $handbackZhDate1
$zh.Range("K2").NumberFormat = $dateFmt
Set-Hyperlink $zh "A2" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md1" $md1
Set-Hyperlink $zh "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6b954ecf0f8160b5e84a6cf6734fb8e3fe98daaa/e2e/$md1" $md1

# Row 3: newly appeared file.
$zh.Range("A3").Value = $md2
$zh.Range("B3").Value = $extDot
$zh.Range("C3").Value = $status
$zh.Range("D3").Value = $extNoDot
$zh.Range("E3").Value = $htStr
$zh.Range("F3").Value = $trueStr
$zh.Range("G3").Value = $zhXlf2
$zh.Range("H3").Value = $handoffDate2
$zh.Range("H3").NumberFormat = $dateFmt
$zh.Range("I3").Value = $md2
$zh.Range("J3").Value = $zhXlf2
$zh.Range("K3").Value = $handbackZhDate2
$zh.Range("K3").NumberFormat = $dateFmt
$zh.Range("L3").Value = ""
$zh.Range("M3").Value = $trueStr
$zh.Range("N3").Value = ""
$zh.Range("O3").Value = $falseStr
$zh.Range("P3").Value = ""
Set-Hyperlink $zh "A3" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md2" $md2
Set-Hyperlink $zh "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/6b954ecf0f8160b5e84a6cf6734fb8e3fe98daaa/e2e/$md2" $md2

$zhTable = $zh.ListObjects.Item(1)
$zhTable.Resize($zh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

# Row 2: re-processed file - new xliff hash + new handback timestamp
# (the "Correspond Handoff Datetime" column reuses the generate date here,
# same as in the source workbook).
$de.Range("A2").Value = $md1
$de.Range("G2").Value = $deXlf1
$de.Range("H2").Value = $genDate1
$de.Range("H2").NumberFormat = $dateFmt
$de.Range("I2").Value = $md1
$de.Range("J2").Value = $deXlf1
$de.Range("K2").Value = $handbackDeDate1
$de.Range("K2").NumberFormat = $dateFmt
Set-Hyperlink $de "A2" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md1" $md1
Set-Hyperlink $de "I2" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f64f95e8895d1904311c3387e0b4219374ff923/e2e/$md1" $md1

# Row 3: newly appeared file.
$de.Range("A3").Value = $md2
$de.Range("B3").Value = $extDot
$de.Range("C3").Value = $status
$de.Range("D3").Value = $extNoDot
$de.Range("E3").Value = $htStr
$de.Range("F3").Value = $trueStr
$de.Range("G3").Value = $deXlf2
$de.Range("H3").Value = $genDate2
$de.Range("H3").NumberFormat = $dateFmt
$de.Range("I3").Value = $md2
$de.Range("J3").Value = $deXlf2
$de.Range("K3").Value = $handbackDeDate2
$de.Range("K3").NumberFormat = $dateFmt
$de.Range("L3").Value = ""
$de.Range("M3").Value = $trueStr
$de.Range("N3").Value = ""
$de.Range("O3").Value = $falseStr
$de.Range("P3").Value = ""
Set-Hyperlink $de "A3" "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/$md2" $md2
Set-Hyperlink $de "I3" "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/9f64f95e8895d1904311c3387e0b4219374ff923/e2e/$md2" $md2

$deTable = $de.ListObjects.Item(1)
$deTable.Resize($de.Range("A1:P3"))
